# Add a new "2022-Q1" sheet before "总计" with fund-holding detail data,
# cloned from the structurally-identical "2021-Q4" sheet so headers,
# column widths/styles and sheet properties match the existing pattern.

$wb = $excel.ActiveWorkbook
$totalSheet = $wb.Worksheets.Item("总计")
$sourceSheet = $wb.Worksheets.Item("2021-Q4")

# Copy() inserts the duplicate immediately before the sheet passed as the
# argument, which is exactly where "2022-Q1" belongs in the tab order.
$sourceSheet.Copy($totalSheet)
$ws = $wb.Worksheets.Item("2021-Q4 (2)")
$ws.Name = "2022-Q1"

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "'166005"
$ws.Range("B2").Style = "Normal"
$ws.Range("C2").Value = "中欧价值发现混合 -A"
$ws.Range("D2").Value = "'43.52"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'93.97"
$ws.Range("E2").Style = "Normal"
$ws.Range("F2").Value = "'3.03"
$ws.Range("F2").Style = "Normal"
$ws.Range("G2").Value = "'1.3187"
$ws.Range("G2").Style = "Normal"
$ws.Range("H2").Value = 9

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "'001882"
$ws.Range("B3").Style = "Normal"
$ws.Range("C3").Value = "中欧价值发现混合 -E"
$ws.Range("D3").Value = "'43.52"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'93.97"
$ws.Range("E3").Style = "Normal"
$ws.Range("F3").Value = "'3.03"
$ws.Range("F3").Style = "Normal"
$ws.Range("G3").Value = "'1.3187"
$ws.Range("G3").Style = "Normal"
$ws.Range("H3").Value = 9

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "'001810"
$ws.Range("B4").Style = "Normal"
$ws.Range("C4").Value = "中欧潜力价值灵活配置混合A"
$ws.Range("D4").Value = "'28.67"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'94.05"
$ws.Range("E4").Style = "Normal"
$ws.Range("F4").Value = "'2.79"
$ws.Range("F4").Style = "Normal"
$ws.Range("G4").Value = "'0.7999"
$ws.Range("G4").Style = "Normal"
$ws.Range("H4").Value = 8

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "'004856"
$ws.Range("B5").Style = "Normal"
$ws.Range("C5").Value = "广发中证全指建筑材料指数A"
$ws.Range("D5").Value = "'13.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'94.61"
$ws.Range("E5").Style = "Normal"
$ws.Range("F5").Value = "'3.16"
$ws.Range("F5").Style = "Normal"
$ws.Range("G5").Value = "'0.4336"
$ws.Range("G5").Style = "Normal"
$ws.Range("H5").Value = 7

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "'004232"
$ws.Range("B6").Style = "Normal"
$ws.Range("C6").Value = "中欧价值发现混合 -C"
$ws.Range("D6").Value = "'10.98"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'93.97"
$ws.Range("E6").Style = "Normal"
$ws.Range("F6").Value = "'3.03"
$ws.Range("F6").Style = "Normal"
$ws.Range("G6").Value = "'0.3327"
$ws.Range("G6").Style = "Normal"
$ws.Range("H6").Value = 9

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "'004857"
$ws.Range("B7").Style = "Normal"
$ws.Range("C7").Value = "广发中证全指建筑材料指数C"
$ws.Range("D7").Value = "'6.05"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'94.61"
$ws.Range("E7").Style = "Normal"
$ws.Range("F7").Value = "'3.16"
$ws.Range("F7").Style = "Normal"
$ws.Range("G7").Value = "'0.1912"
$ws.Range("G7").Style = "Normal"
$ws.Range("H7").Value = 7

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "'166024"
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Value = "中欧恒利三年定期开放混合"
$ws.Range("D8").Value = "'4.48"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'98.71"
$ws.Range("E8").Style = "Normal"
$ws.Range("F8").Value = "'3.73"
$ws.Range("F8").Style = "Normal"
$ws.Range("G8").Value = "'0.1671"
$ws.Range("G8").Style = "Normal"
$ws.Range("H8").Value = 7

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "'159745"
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").Value = "国泰中证全指建筑材料交易型开放式指数证券投资基金"
$ws.Range("D9").Value = "'3.76"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'98.37"
$ws.Range("E9").Style = "Normal"
$ws.Range("F9").Value = "'3.39"
$ws.Range("F9").Style = "Normal"
$ws.Range("G9").Value = "'0.1275"
$ws.Range("G9").Style = "Normal"
$ws.Range("H9").Value = 7

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "'001050"
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Value = "汇添富成长多因子量化策略股票"
$ws.Range("D10").Value = "'11.48"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'92.68"
$ws.Range("E10").Style = "Normal"
$ws.Range("F10").Value = "'1.01"
$ws.Range("F10").Style = "Normal"
$ws.Range("G10").Value = "'0.1159"
$ws.Range("G10").Style = "Normal"
$ws.Range("H10").Value = 2

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "'005764"
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = "中欧潜力价值灵活配置混合C"
$ws.Range("D11").Value = "'3.43"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'94.05"
$ws.Range("E11").Style = "Normal"
$ws.Range("F11").Value = "'2.79"
$ws.Range("F11").Style = "Normal"
$ws.Range("G11").Value = "'0.0957"
$ws.Range("G11").Style = "Normal"
$ws.Range("H11").Value = 8

$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "'001891"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = "中欧成长优选回报灵活配置混合E"
$ws.Range("D12").Value = "'2.97"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'94.42"
$ws.Range("E12").Style = "Normal"
$ws.Range("F12").Value = "'2.67"
$ws.Range("F12").Style = "Normal"
$ws.Range("G12").Value = "'0.0793"
$ws.Range("G12").Style = "Normal"
$ws.Range("H12").Value = 8

$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "'166020"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = "中欧成长优选回报灵活配置混合A"
$ws.Range("D13").Value = "'2.97"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'94.42"
$ws.Range("E13").Style = "Normal"
$ws.Range("F13").Value = "'2.67"
$ws.Range("F13").Style = "Normal"
$ws.Range("G13").Value = "'0.0793"
$ws.Range("G13").Style = "Normal"
$ws.Range("H13").Value = 8

$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "'516750"
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = "富国中证全指建筑材料ETF"
$ws.Range("D14").Value = "'0.47"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'98.22"
$ws.Range("E14").Style = "Normal"
$ws.Range("F14").Value = "'3.52"
$ws.Range("F14").Style = "Normal"
$ws.Range("G14").Value = "'0.0165"
$ws.Range("G14").Style = "Normal"
$ws.Range("H14").Value = 6

$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "'004135"
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Value = "申万菱信量化成长混合"
$ws.Range("D15").Value = "'0.49"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'86.91"
$ws.Range("E15").Style = "Normal"
$ws.Range("F15").Value = "'1.98"
$ws.Range("F15").Style = "Normal"
$ws.Range("G15").Value = "'0.0097"
$ws.Range("G15").Style = "Normal"
$ws.Range("H15").Value = 6

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "'164811"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = "工银瑞信中证京津冀协同发展主题指数（LOF）A"
$ws.Range("D16").Value = "'0.23"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'94.28"
$ws.Range("E16").Style = "Normal"
$ws.Range("F16").Value = "'2.89"
$ws.Range("F16").Style = "Normal"
$ws.Range("G16").Value = "'0.0066"
$ws.Range("G16").Style = "Normal"
$ws.Range("H16").Value = 10

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "'164825"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "工银瑞信中证京津冀协同发展主题指数（LOF）C"
$ws.Range("D17").Value = "'0.06"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'94.28"
$ws.Range("E17").Style = "Normal"
$ws.Range("F17").Value = "'2.89"
$ws.Range("F17").Style = "Normal"
$ws.Range("G17").Value = "'0.0017"
$ws.Range("G17").Style = "Normal"
$ws.Range("H17").Value = 10
